$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set F7 to the new string value "done" (adds a shared string and updates the cell)
$ws.Range("F7").Value = "done"

# Update the selected cell to G7 (to match the new selection in the sheet view)
$ws.Range("G7").Select()
